$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "39+37=" "88-63="
Replace-Text "46+3=" "16+15="
Replace-Text "4-2=" "43-2="
Replace-Text "93-63=" "72+9="
Replace-Text "15-11=" "97-72="
Replace-Text "65-49=" "97-30="
Replace-Text "36-15=" "96-13="
Replace-Text "73-50=" "29+64="
Replace-Text "67+3=" "5-1="
Replace-Text "20+71=" "78-21="
Replace-Text "31+23=" "90-86="
Replace-Text "4+1=" "92-67="
Replace-Text "41+5=" "76+17="
Replace-Text "80-44=" "49+18="
Replace-Text "15+63=" "68+28="
Replace-Text "64-19=" "68-22="
Replace-Text "49-3=" "25+20="
Replace-Text "71+10=" "12+5="
Replace-Text "63-0=" "9+80="
Replace-Text "54+10=" "70+21="
Replace-Text "41+42=" "70-15="
Replace-Text "89-52=" "89-73="
Replace-Text "90-1=" "79-50="
Replace-Text "57+36=" "28+63="
Replace-Text "34+44=" "39+23="
Replace-Text "51+23=" "46-45="
Replace-Text "34-21=" "25+27="
Replace-Text "0+16=" "84-37="
Replace-Text "24+41=" "68+12="
Replace-Text "83-47=" "30-26="
Replace-Text "75+5=" "47+52="
Replace-Text "74-24=" "80-38="
Replace-Text "92-57=" "71-70="
Replace-Text "83-67=" "84-5="
Replace-Text "28-3=" "27+5="
Replace-Text "98-48=" "28+14="
Replace-Text "44-11=" "1+71="
Replace-Text "84+5=" "73+4="
Replace-Text "18-5=" "93-25="
Replace-Text "5+14=" "86-41="
Replace-Text "37-3=" "66-28="
Replace-Text "22+45=" "72-24="
Replace-Text "44+9=" "94-40="
Replace-Text "54+6=" "41-1="
Replace-Text "53+31=" "10+29="
Replace-Text "71+27=" "94-84="
Replace-Text "56+23=" "76-45="
Replace-Text "15+84=" "87-75="
Replace-Text "63-38=" "58-54="
Replace-Text "71-1=" "77-76="
Replace-Text "13+67=" "61+26="
Replace-Text "63-57=" "62+7="
Replace-Text "21+61=" "31+41="
Replace-Text "60-37=" "2+47="
Replace-Text "6+70=" "61-14="
Replace-Text "63+35=" "21+65="
Replace-Text "39+27=" "53-35="
Replace-Text "20-14=" "36-10="
Replace-Text "71-20=" "1+63="
Replace-Text "34+45=" "0+9="
Replace-Text "48+0=" "23-22="
Replace-Text "21+62=" "66-6="
Replace-Text "88+6=" "6+91="
Replace-Text "0+54=" "64-4="
Replace-Text "1+11=" "43-23="
Replace-Text "65-26=" "40+44="
Replace-Text "57-47=" "63-17="
Replace-Text "67-15=" "36+24="
Replace-Text "88-59=" "35+18="
Replace-Text "16+31=" "7+51="
Replace-Text "10+71=" "12+8="
Replace-Text "71-53=" "98-97="
Replace-Text "99-42=" "32+21="
Replace-Text "89+2=" "79-26="
Replace-Text "14+2=" "43-34="
Replace-Text "80-64=" "99-1="
Replace-Text "68-7=" "94-48="
Replace-Text "56-47=" "19+70="
Replace-Text "16+75=" "39+60="
Replace-Text "57-22=" "88-58="
Replace-Text "74+9=" "36-20="
Replace-Text "66+15=" "53-19="
Replace-Text "3+70=" "48+7="
Replace-Text "8+16=" "69-32="
Replace-Text "83-63=" "98-0="
Replace-Text "34+5=" "76-30="
Replace-Text "8+24=" "76+13="
Replace-Text "56-38=" "3+8="
Replace-Text "49+22=" "12+25="
Replace-Text "86-6=" "43+3="
Replace-Text "87-73=" "39+56="
Replace-Text "65-9=" "2+78="
Replace-Text "43-36=" "48-2="
Replace-Text "83-20=" "45-14="
Replace-Text "14+48=" "21+76="
Replace-Text "58-1=" "86+3="
Replace-Text "35-28=" "56+36="
Replace-Text "40+2=" "57-26="
Replace-Text "99-38=" "35+2="
Replace-Text "88-1=" "44+39="
